$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Delete the "Senses of Scripture" section (an empty ListParagraph,
#    the "Senses of Scripture" heading, and its four sub-bullets), which
#    were trimmed from the end of the outline, right after the
#    "Usefulness of Scripture" bullet.
# ---------------------------------------------------------------------
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Usefulness of Scripture*") {
        $startPara = $i + 1
    }
    if ($t -like "*eternal significance*") {
        $endPara = $i
    }
}
if ($startPara -ne $null -and $endPara -ne $null -and $endPara -ge $startPara) {
    $delRange = $d.Range($d.Paragraphs($startPara).Range.Start, $d.Paragraphs($endPara).Range.End)
    $delRange.Delete()
}

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" last-edit-location bookmark from inside
#    "courageous" to inside "Utilitarianism" (between "Utilitariani"
#    and "sm"), reflecting where the author's final edit landed.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Utilitarianism")
if ($rng.Find.Found) {
    $splitPos = $rng.Start + 12
    $d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))
}

# ---------------------------------------------------------------------
# 3) The footer's cached " PAGE " field result drops from "4" to "1"
#    now that the outline is shorter. Update the cached field result
#    text directly (this simulator has no live pagination engine to
#    recompute it for us).
# ---------------------------------------------------------------------
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
foreach ($fld in $ftr.Range.Fields) {
    if ($fld.Type -eq 33 -or $fld.Code.Text -like "*PAGE*") {
        $resultRange = $fld.Result
        if ($resultRange.Characters.Count -ge 1) {
            $resultRange.Characters(1).Text = "1"
        }
    }
}
